# Upload excel files with prices
# Swap the two rows (19 and 20) that got re-ordered in the scraped data,
# then refresh the scrape timestamp (column O) for every data row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Swap contents of row 19 and row 20 (columns A:N) ---
# Columns E and F hold true numbers; every other column in this sheet is
# stored as text (even when the text looks numeric, e.g. ids/prices), so
# force those columns to keep a text number format while we shuffle values.
$colCount = 14  # columns A..N
$textColumns = @(1, 2, 3, 4, 7, 8, 9, 10, 11, 12, 13, 14)  # A,B,C,D,G,H,I,J,K,L,M,N

foreach ($c in $textColumns) {
    $ws.Cells.Item(19, $c).NumberFormat = "@"
    $ws.Cells.Item(20, $c).NumberFormat = "@"
}

for ($c = 1; $c -le $colCount; $c++) {
    $cell19 = $ws.Cells.Item(19, $c)
    $cell20 = $ws.Cells.Item(20, $c)
    $tmp = $cell19.Value2
    $cell19.Value2 = $cell20.Value2
    $cell20.Value2 = $tmp
}

# --- Update timestamp column (O) for all data rows (2 through 33) ---
$newTimestamp = "2022-08-18 20:58:33"
for ($r = 2; $r -le 33; $r++) {
    $ws.Cells.Item($r, 15).Value2 = $newTimestamp
}
